## BATT_DCDC.xlsx — version 5
## "se añaden los componentes correspondientes a los conectores interplaca."
##
## 1) BOM sheet: NEW flag for RSENSE (row 30 / B36) changes from "SI" to "YES"
## 2) BOM sheet: add position 38 -> J1 connector (DF40C-20DP-0.4V_51_)
## 3) _HISTORY sheet: fix the "5-ene-2023" text date into a real date value,
##    and log a new history entry for this change (version 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BATT_DCDC")
$hist = $wb.Worksheets.Item("_HISTORY")

# --- 1) Normalize "SI" -> "YES" on the RSENSE row --------------------------
$ws.Range("B36").Value = "YES"

# --- 2) Add the new BOM row (position 38 / J1 connector) -------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null

$ws.Range("A44").Value = 38
$ws.Range("B44").Value = "YES"
$ws.Range("B44").HorizontalAlignment = -4108
$ws.Range("C44").Value = 1
$ws.Range("D44").Value = "J1"
$ws.Range("E44").Value = "DF40C-20DP-0.4V_51_"
$ws.Range("F44").Value = "CONNECTOR"
$ws.Range("G44").Value = "DF40C-20DPA"
$ws.Range("J44").Value = "Se añade la posicion 38."
$ws.Hyperlinks.Add($ws.Range("J44"), "https://www.digikey.es/es/products/detail/hirose-electric-co-ltd/DF40C-20DP-0-4V-51/1969479")
$ws.Range("J44").Style = $ws.Range("J38").Style

# --- extend the NO/YES conditional formatting over the new row -------------
$oldRng = $ws.Range("B7:B43")
$oldRng.FormatConditions.Delete()
$newRng = $ws.Range("B7:B44")

$fc1 = $newRng.FormatConditions.Add(9, 0, "NO")
$fc1.Formula1 = 'NOT(ISERROR(SEARCH("NO",B7)))'
$fc1.Text = "NO"

$fc2 = $newRng.FormatConditions.Add(9, 0, "YES")
$fc2.Formula1 = 'NOT(ISERROR(SEARCH("YES",B7)))'
$fc2.Text = "YES"

$fc3 = $newRng.FormatConditions.Add(9, 0, "NO")
$fc3.Formula1 = 'NOT(ISERROR(SEARCH("NO",B7)))'
$fc3.Text = "NO"

# --- 3) _HISTORY sheet updates ---------------------------------------------
# row 7 ("5-ene-2023" text) becomes a real date, matching rows 4-6
$hist.Range("B7").Value = 44931
$hist.Range("B7").NumberFormat = "d-mmm"

# new history row: version 5
$hist.Range("A8").Value = 5
$hist.Range("B8").Value = 44958
$hist.Range("B8").NumberFormat = "d-mmm"
$hist.Range("C8").Value = "JRC"
$hist.Range("D8").Value = "Se añade la posicion 38."
